# Commit: "remove links from workbooks to FixedIncome.xla"
#
# The workbook has one external link (an externalReference in
# xl/workbook.xml backed by xl/externalLinks/externalLink1.xml, pointing
# at .../framework/addin/Menu.xla) that is used by a single formula:
#   'General Settings'!D8  =[1]!qlSerializationPath(Trigger)
# which currently resolves to "C:\Projects\quantlib\...\010_Quotes\".
#
# We replace that formula with the static (local) value the author baked
# in ("C:\Users\erik\junk\") and then sever the external link entirely so
# the externalReferences/externalLinks parts disappear from the package.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General Settings")
$cell = $ws.Range("D8")

# Leading apostrophe = Excel's "treat as text" marker. It forces the cell
# to keep its existing quotePrefix'd style (rather than Excel silently
# switching the cell to the no-quote-prefix twin style it otherwise
# reassigns when a formula cell is overwritten with a literal string) and
# is stripped from the stored text itself.
$cell.Value = "'C:\Users\erik\junk\"

# Break every external link/source so the externalReferences node and the
# externalLinks part are removed from the saved package.
$sources = $wb.LinkSources(1)
if ($sources) {
    foreach ($source in $sources) {
        $wb.BreakLink($source, 1)
    }
}
